$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 3856
$ws.Range("E2").Value = 128
$ws.Range("F2").Value = 128
$ws.Range("G2").Value = 132
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 3984
$ws.Range("L2").Value = 2731
$ws.Range("M2").Value = 1254
$ws.Range("N2").Value = 1254
$ws.Range("P2").Value = 176
$ws.Range("Q2").Value = 67
$ws.Range("R2").Value = 114
$ws.Range("S2").Value = -26
$ws.Range("T2").Value = 14
$ws.Range("U2").Value = 52
$ws.Range("V2").Value = 38
$ws.Range("W2").Value = 3.33
$ws.Range("X2").Value = 2.47
$ws.Range("Y2").Value = 7.85
$ws.Range("Z2").Value = 2.36
$ws.Range("AA2").Value = 217.78
$ws.Range("AB2").Value = 667.87
$ws.Range("AC2").Value = 584
$ws.Range("AD2").Value = 13.17
$ws.Range("AE2").Value = 7736
$ws.Range("AF2").Value = 0.99
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.6
$ws.Range("AI2").Value = 33.52
$ws.Range("AJ2").Value = 16567409
$ws.Range("O2").ClearContents()

$ws.Range("D3").Value = 4793
$ws.Range("E3").Value = 138
$ws.Range("F3").Value = 138
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = 116
$ws.Range("I3").Value = 116
$ws.Range("K3").Value = 4828
$ws.Range("L3").Value = 3479
$ws.Range("M3").Value = 1349
$ws.Range("N3").Value = 1349
$ws.Range("P3").Value = 176
$ws.Range("Q3").Value = 179
$ws.Range("R3").Value = -55
$ws.Range("S3").Value = -29
$ws.Range("T3").Value = 22
$ws.Range("U3").Value = 157
$ws.Range("V3").Value = 106
$ws.Range("W3").Value = 2.88
$ws.Range("X3").Value = 2.42
$ws.Range("Y3").Value = 8.91
$ws.Range("Z3").Value = 2.63
$ws.Range("AA3").Value = 257.86
$ws.Range("AB3").Value = 719.89
$ws.Range("AC3").Value = 700
$ws.Range("AD3").Value = 12.57
$ws.Range("AE3").Value = 8325
$ws.Range("AF3").Value = 1.06
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 2.27
$ws.Range("AI3").Value = 27.94
$ws.Range("AJ3").Value = 16567409
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

$ws.Range("D4").Value = 4080
$ws.Range("E4").Value = 196
$ws.Range("F4").Value = 196
$ws.Range("G4").Value = 207
$ws.Range("H4").Value = 149
$ws.Range("I4").Value = 149
$ws.Range("K4").Value = 5071
$ws.Range("L4").Value = 3610
$ws.Range("M4").Value = 1461
$ws.Range("N4").Value = 1461
$ws.Range("P4").Value = 176
$ws.Range("Q4").Value = 136
$ws.Range("R4").Value = -140
$ws.Range("S4").Value = -38
$ws.Range("T4").Value = 17
$ws.Range("U4").Value = 119
$ws.Range("V4").Value = 47
$ws.Range("W4").Value = 4.8
$ws.Range("X4").Value = 3.66
$ws.Range("Y4").Value = 10.64
$ws.Range("Z4").Value = 3.02
$ws.Range("AA4").Value = 247.12
$ws.Range("AB4").Value = 786.88
$ws.Range("AC4").Value = 902
$ws.Range("AD4").Value = 9.06
$ws.Range("AE4").Value = 9013
$ws.Range("AF4").Value = 0.91
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 3.06
$ws.Range("AI4").Value = 27.11
$ws.Range("AJ4").Value = 16567409
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

$ws.Range("D5").Value = 5013
$ws.Range("E5").Value = 197
$ws.Range("F5").Value = 197
$ws.Range("G5").Value = 165
$ws.Range("H5").Value = 117
$ws.Range("I5").Value = 117
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6037
$ws.Range("L5").Value = 4818
$ws.Range("M5").Value = 1219
$ws.Range("N5").Value = 1219
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 176
$ws.Range("Q5").Value = 515
$ws.Range("R5").Value = -152
$ws.Range("S5").Value = -36
$ws.Range("T5").Value = 26
$ws.Range("U5").Value = 490
$ws.Range("V5").Value = 47
$ws.Range("W5").Value = 3.93
$ws.Range("X5").Value = 2.33
$ws.Range("Y5").Value = 8.71
$ws.Range("Z5").Value = 2.1
$ws.Range("AA5").Value = 395.31
$ws.Range("AB5").Value = 652.38
$ws.Range("AC5").Value = 705
$ws.Range("AD5").Value = 14.62
$ws.Range("AE5").Value = 7519
$ws.Range("AF5").Value = 1.37
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 2.43
$ws.Range("AI5").Value = 34.71
$ws.Range("AJ5").Value = 16567409

$ws.Range("D6").Value = 5382
$ws.Range("E6").Value = 254
$ws.Range("F6").Value = 254
$ws.Range("G6").Value = 226
$ws.Range("H6").Value = 153
$ws.Range("I6").Value = 153
$ws.Range("K6").Value = 6328
$ws.Range("L6").Value = 5024
$ws.Range("M6").Value = 1304
$ws.Range("N6").Value = 1304
$ws.Range("P6").Value = 176
$ws.Range("Q6").Value = -420
$ws.Range("R6").Value = 17
$ws.Range("S6").Value = 12
$ws.Range("T6").Value = 26
$ws.Range("U6").Value = -446
$ws.Range("V6").Value = 99
$ws.Range("W6").Value = 4.72
$ws.Range("X6").Value = 2.84
$ws.Range("Y6").Value = 12.11
$ws.Range("Z6").Value = 2.47
$ws.Range("AA6").Value = 385.21
$ws.Range("AB6").Value = 702.3
$ws.Range("AC6").Value = 922
$ws.Range("AD6").Value = 8.16
$ws.Range("AE6").Value = 8046
$ws.Range("AF6").Value = 0.94
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 3.32
$ws.Range("AI6").Value = 26.52
$ws.Range("AJ6").Value = 16567409

$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

